{"js": "// Apply the \"football-glory\" copy edits: refreshed title/meta text and\n// expanded \"what we like\" / \"what we don't like\" bullet points.\nconst replacements = [\n  {\n    from: \"Play Football Glory for Free - Unique Gameplay and Exciting Bonus Features\",\n    to: \"Play Football Glory - Free Online Slot Game\",\n  },\n  {\n    from: \"Well-crafted graphics\",\n    to: \"Well-crafted graphics with realistic illustrations and vivid colors\",\n  },\n  {\n    from: \"Unique game mechanics\",\n    to: \"Unique 5x4 layout that sets it apart from other online slots\",\n  },\n  {\n    from: \"Special symbols for bigger wins\",\n    to: \"Special symbols like Wild and Free Spins for bigger wins\",\n  },\n  {\n    from: \"Exciting Free Spins bonus mode\",\n    to: \"Exciting bonus features, including the Free Spins mode\",\n  },\n  {\n    from: \"Limited football theme\",\n    to: \"Limited selection of football-inspired slots\",\n  },\n  {\n    from: \"No jackpot feature\",\n    to: \"May not appeal to non-football fans\",\n  },\n  {\n    from:\n      \"Read our review of Football Glory, the online slot game with well-crafted graphics, unique mechanics, and exciting bonus features. Play for free now!\",\n    to:\n      \"Read our review of Football Glory, a free online slot game with exciting bonus features.\",\n  },\n];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"football-glory\" copy edits: refreshed title/meta text and\n# expanded \"what we like\" / \"what we don't like\" bullet points.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.Execute(\n        $findText,     # FindText\n        $true,         # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap (wdFindContinue)\n        $false,        # Format\n        $replaceText,  # ReplaceWith\n        2              # Replace (wdReplaceAll)\n    )\n}\n\nReplace-Text \"Play Football Glory for Free - Unique Gameplay and Exciting Bonus Features\" \"Play Football Glory - Free Online Slot Game\"\nReplace-Text \"Well-crafted graphics\" \"Well-crafted graphics with realistic illustrations and vivid colors\"\nReplace-Text \"Unique game mechanics\" \"Unique 5x4 layout that sets it apart from other online slots\"\nReplace-Text \"Special symbols for bigger wins\" \"Special symbols like Wild and Free Spins for bigger wins\"\nReplace-Text \"Exciting Free Spins bonus mode\" \"Exciting bonus features, including the Free Spins mode\"\nReplace-Text \"Limited football theme\" \"Limited selection of football-inspired slots\"\nReplace-Text \"No jackpot feature\" \"May not appeal to non-football fans\"\nReplace-Text \"Read our review of Football Glory, the online slot game with well-crafted graphics, unique mechanics, and exciting bonus features. Play for free now!\" \"Read our review of Football Glory, a free online slot game with exciting bonus features.\"\n"}
